$wb = $excel.ActiveWorkbook

# Map of sheet name -> new price value for the appended 2025-05-01 row.
$updates = @{
    "N-Dense"                   = "38"
    "N-Type"                    = "37.3"
    "N-type Wafer"              = "1.02"
    "Cell Topcon 183mm"         = "0.273"
    "Module Topcon 183mm"       = "0.09"
    "Silver Rear_side"          = "5,360"
    "Silver Busbar front-side"  = "8,025"
    "Silver finger front-side"  = "8,075"
    "USD_CNY"                   = "7.2927"
}

foreach ($name in $updates.Keys) {
    $ws = $wb.Worksheets.Item($name)
    $dateCell = $ws.Cells.Item(61, 1)
    $valueCell = $ws.Cells.Item(61, 2)

    # Force text entry (leading apostrophe) so cells are stored the same way
    # as the rest of the Date/Price columns (inline/shared string, not a
    # date-serial number or General number), then strip the quote-prefix
    # style that forcing text entry leaves behind so no extra cell style is
    # introduced.
    $dateCell.Value = "'2025-05-01"
    $valueCell.Value = "'" + $updates[$name]

    $dateCell.Style = "Normal"
    $valueCell.Style = "Normal"
}
